$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1353.4
$ws.Range("J17").Value = 1353.4
$ws.Range("L17").Value = 4060.2
$ws.Range("N17").Value = -4396.200000000001
$ws.Range("H32").Value = 641.8570999999999
$ws.Range("I32").Value = 400
$ws.Range("J32").Value = 738.6
$ws.Range("K32").Value = 400
$ws.Range("L32").Value = 738.6
$ws.Range("M32").Value = -74
$ws.Range("N32").Value = -1390.6
$ws.Range("H40").Value = 2750.1667
$ws.Range("I40").Value = 1785.7142
$ws.Range("J40").Value = 4100.4
$ws.Range("K40").Value = 1785.7142
$ws.Range("L40").Value = 4100.4
$ws.Range("M40").Value = -1610.7142
$ws.Range("N40").Value = -4450.4
$ws.Range("H51").Value = 6627.92
$ws.Range("I51").Value = 12020
$ws.Range("J51").Value = 3033.2
$ws.Range("K51").Value = 12020
$ws.Range("L51").Value = 3033.2
$ws.Range("M51").Value = -11536
$ws.Range("N51").Value = -4001.2
$ws.Range("H64").Value = 44587.418
$ws.Range("I64").Value = 69379.87
$ws.Range("J64").Value = 3266.6667
$ws.Range("K64").Value = 69379.87
$ws.Range("L64").Value = 3266.6667
$ws.Range("M64").Value = -69131.87
$ws.Range("N64").Value = -3762.6667
$ws.Range("H67").Value = 44587.418
$ws.Range("I67").Value = 69379.87
$ws.Range("J67").Value = 3266.6667
$ws.Range("K67").Value = 69379.87
$ws.Range("L67").Value = 3266.6667
$ws.Range("M67").Value = -68521.87
$ws.Range("N67").Value = -4982.6667
$ws.Range("H127").Value = 26317784
$ws.Range("I127").Value = 549
$ws.Range("J127").Value = 29413930
$ws.Range("K127").Value = 1647
$ws.Range("L127").Value = 88241790
$ws.Range("M127").Value = 3313
$ws.Range("N127").Value = -88251710
$ws.Range("H131").Value = 4304.163
$ws.Range("I131").Value = 359
$ws.Range("J131").Value = 4752.477
$ws.Range("K131").Value = 1077
$ws.Range("L131").Value = 14257.431
$ws.Range("M131").Value = 3963
$ws.Range("N131").Value = -24337.431
$ws.Range("H135").Value = 1362.6
$ws.Range("I135").Value = 370.86957
$ws.Range("J135").Value = 2704.353
$ws.Range("K135").Value = 3337.82613
$ws.Range("L135").Value = 24339.177
$ws.Range("M135").Value = -802.8261299999999
$ws.Range("N135").Value = -29409.177
$ws.Range("H141").Value = 4360.909
$ws.Range("I141").Value = 4849.1665
$ws.Range("J141").Value = 3775
$ws.Range("K141").Value = 14547.4995
$ws.Range("L141").Value = 11325
$ws.Range("M141").Value = -9367.499500000002
$ws.Range("N141").Value = -21685

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2233.5757
$ws.Range("J61").Value = 2986.6667
$ws.Range("L61").Value = 2986.6667
$ws.Range("N61").Value = -3410.6667
$ws.Range("H136").Value = 2233.5757
$ws.Range("J136").Value = 2986.6667
$ws.Range("L136").Value = 8960.000100000001
$ws.Range("N136").Value = -14060.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 31500
$ws.Range("J112").Value = 31500
$ws.Range("L112").Value = 31500
$ws.Range("N112").Value = -34454
$ws.Range("H134").Value = 2308.838
$ws.Range("I134").Value = 2052.4666
$ws.Range("J134").Value = 3407.5715
$ws.Range("K134").Value = 6157.399800000001
$ws.Range("L134").Value = 10222.7145
$ws.Range("M134").Value = -3622.399800000001
$ws.Range("N134").Value = -15292.7145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2913.76
$ws.Range("J31").Value = 3715.4443
$ws.Range("L31").Value = 3715.4443
$ws.Range("N31").Value = -4305.4443
$ws.Range("H34").Value = 2913.76
$ws.Range("J34").Value = 3715.4443
$ws.Range("L34").Value = 3715.4443
$ws.Range("N34").Value = -4119.4443
$ws.Range("H58").Value = 2541.4375
$ws.Range("I58").Value = 2669.9
$ws.Range("J58").Value = 2327.3333
$ws.Range("K58").Value = 2669.9
$ws.Range("L58").Value = 2327.3333
$ws.Range("M58").Value = -2466.9
$ws.Range("N58").Value = -2733.3333
$ws.Range("H99").Value = 2281.5
$ws.Range("I99").Value = 1698.4
$ws.Range("J99").Value = 2605.4443
$ws.Range("K99").Value = 1698.4
$ws.Range("L99").Value = 2605.4443
$ws.Range("M99").Value = -200.4000000000001
$ws.Range("N99").Value = -5601.4443
$ws.Range("H126").Value = 2281.5
$ws.Range("I126").Value = 1698.4
$ws.Range("J126").Value = 2605.4443
$ws.Range("K126").Value = 5095.200000000001
$ws.Range("L126").Value = 7816.3329
$ws.Range("M126").Value = -2625.200000000001
$ws.Range("N126").Value = -12756.3329
$ws.Range("H134").Value = 1339.9286
$ws.Range("I134").Value = 887.1818
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 2661.5454
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -126.5454
$ws.Range("N134").Value = -14070
$ws.Range("H136").Value = 2541.4375
$ws.Range("I136").Value = 2669.9
$ws.Range("J136").Value = 2327.3333
$ws.Range("K136").Value = 8009.700000000001
$ws.Range("L136").Value = 6981.999899999999
$ws.Range("M136").Value = -5459.700000000001
$ws.Range("N136").Value = -12081.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 81.375
$ws.Range("I38").Value = 80
$ws.Range("K38").Value = 240
$ws.Range("M38").Value = 107
$ws.Range("H69").Value = 1633.3334
$ws.Range("I69").Value = 800
$ws.Range("K69").Value = 2400
$ws.Range("M69").Value = -1589
$ws.Range("H72").Value = 1633.3334
$ws.Range("I72").Value = 800
$ws.Range("K72").Value = 7200
$ws.Range("M72").Value = -3144

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1219.9445
$ws.Range("I22").Value = 1893.1666
$ws.Range("K22").Value = 1893.1666
$ws.Range("M22").Value = -1598.1666
$ws.Range("H27").Value = 1219.9445
$ws.Range("I27").Value = 1893.1666
$ws.Range("K27").Value = 1893.1666
$ws.Range("M27").Value = -1786.1666
$ws.Range("H136").Value = 1668.6285
$ws.Range("J136").Value = 2129.9167
$ws.Range("L136").Value = 6389.750100000001
$ws.Range("N136").Value = -11489.7501

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2052.2068
$ws.Range("I132").Value = 1895.8334
$ws.Range("J132").Value = 2802.8
$ws.Range("K132").Value = 5687.5002
$ws.Range("L132").Value = 8408.400000000001
$ws.Range("M132").Value = -3157.5002
$ws.Range("N132").Value = -13468.4
